$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the 5th employee record (row 6) to the "Reporte_Empleados" sheet.

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "carlos "
$ws.Range("C6").Value = "luis"
$ws.Range("D6").Value = "Alvarado "
$ws.Range("E6").Value = "Ken"

# F6, I6, O6 and T6 hold digit-only / date-shaped text that Excel would
# otherwise auto-convert to a number or a date. Force text formatting
# first, then strip the resulting style back to Normal so the cell keeps
# a plain string value with no style index, matching the rest of the
# data rows.
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "12435678"
$ws.Range("F6").Style = "Normal"

$ws.Range("G6").Value = "fefsfsfes"
$ws.Range("H6").Value = "fsfsfsf"

$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1999-11-11"
$ws.Range("I6").Style = "Normal"

$ws.Range("J6").Value = 25
$ws.Range("K6").Value = "Masculino"
$ws.Range("L6").Value = "Casado"
$ws.Range("M6").Value = "Analista"
$ws.Range("N6").Value = "Comercialización"

$ws.Range("O6").NumberFormat = "@"
$ws.Range("O6").Value = "12-11-2023"
$ws.Range("O6").Style = "Normal"

$ws.Range("P6").Value = "Principal"
$ws.Range("Q6").Value = "Quincenal"
$ws.Range("R6").Value = "Activo"
$ws.Range("S6").Value = "Banplus"

$ws.Range("T6").NumberFormat = "@"
$ws.Range("T6").Value = "31231414141414"
$ws.Range("T6").Style = "Normal"
